$wb = $excel.ActiveWorkbook

# Remember which sheet is active so we can restore the selection/tab at the end.
$originalActiveName = $wb.ActiveSheet.Name

# ---------------------------------------------------------------------------
# Step 1: insert a new "2022-Q1" worksheet right before the "总计" sheet.
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($zongji)
$newSheet.Name = "2022-Q1"

# Match the page setup used by the other per-quarter/summary sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Header row (B1:H1), copying the bold/bordered look from the template sheet.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund holdings data for 2022-Q1.
$fundRows = @(
    @("012124", "博道盛彦混合型证券投资基金A", "3.40", "88.74", "2.95", "0.1003", 10),
    @("011685", "创金合信先进装备股票A", "0.73", "92.01", "9.49", "0.0693", 2),
    @("257050", "国联安主题驱动混合", "1.50", "65.37", "3.32", "0.0498", 7),
    @("006803", "嘉实互通精选股票", "0.61", "89.59", "4.72", "0.0288", 10),
    @("011686", "创金合信先进装备股票C", "0.17", "92.01", "9.49", "0.0161", 2),
    @("012125", "博道盛彦混合型证券投资基金C", "0.15", "88.74", "2.95", "0.0044", 10)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    # Index column (A) - numeric, styled like the template's index column.
    $template.Range("A2").Copy()
    $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $newSheet.Cells.Item($r, 1).Value = $i

    # Text columns (B..G) must stay text, matching the source data.
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Style = "Normal"

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    # Rank column (H) - numeric.
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with a new 2022-Q1 row, pushing the
# existing rows down (and renumbering the leading index column).
#
# NOTE: `$zongji` was captured *before* the insert above and was only used to
# position the new sheet (`Worksheets.Add(Before:=...)`); worksheet
# references are positional, so after a sheet gets inserted at that same
# slot, `$zongji` no longer points at "总计" itself. Re-resolve it by name.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 4; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Cells.Item($dst, 2).Value = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($dst, 3).Value = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($dst, 4).Value = $summary.Cells.Item($r, 4).Value2
}

# Give the newly created last row (A5) the same index-column styling as the rest.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)

for ($r = 2; $r -le 5; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 0.27

# ---------------------------------------------------------------------------
# Step 3: restore the originally active sheet/selection.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($originalActiveName).Activate()

Write-Output "done"
